$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Price" (column D) values, forcing them to remain text ---
# (many look like plain numbers/dates to Excel, so we set an explicit
#  text format before assigning, then restore the default "Normal" style)
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '22.101.15'
$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.557.93'
$ws.Range('D3').Style = 'Normal'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('D4').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '292.47'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3989'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3236'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '44.03'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07324'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.711'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.00001141'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.659'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.554.11'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06602'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '83.88'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9998'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.316'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '15.75'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.31'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '22.109.69'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.359'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.452'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '148.72'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.69'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.865'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.729.07'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '119.16'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.013'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.750'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08376'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.624'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '9.103'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02277'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.06140'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.145'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.220'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.0000'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '10.79'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.5872'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.768'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '13.06'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '118.99'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.143'
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06854'
$ws.Range('D51').Style = 'Normal'

# --- Update Coin / Link / Volume(1h) columns (plain text, no conversion risk) ---
$ws.Range('E2').Value = '  -0.86%  '
$ws.Range('E3').Value = '  -0.21%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('E5').Value = '  +0.09%  '
$ws.Range('E6').Value = '  +1.27%  '
$ws.Range('E7').Value = '  +5.64%  '
$ws.Range('E8').Value = '  -1.51%  '
$ws.Range('E9').Value = '  -1.59%  '
$ws.Range('E10').Value = '  -1.09%  '
$ws.Range('E11').Value = '  -5.55%  '
$ws.Range('E12').Value = '  +0.00%  '
$ws.Range('E13').Value = '  -7.19%  '
$ws.Range('E14').Value = '  -2.89%  '
$ws.Range('E15').Value = '  +5.67%  '
$ws.Range('E16').Value = '  -1.81%  '
$ws.Range('E17').Value = '  +0.30%  '
$ws.Range('E18').Value = '  -0.85%  '
$ws.Range('E19').Value = '  -2.96%  '
$ws.Range('E20').Value = '  +0.02%  '
$ws.Range('E21').Value = '  -1.73%  '
$ws.Range('E22').Value = '  -2.94%  '
$ws.Range('E23').Value = '  -3.61%  '
$ws.Range('E24').Value = '  -0.76%  '
$ws.Range('E25').Value = '  +2.71%  '
$ws.Range('E26').Value = '  -6.32%  '
$ws.Range('E27').Value = '  -1.88%  '
$ws.Range('E28').Value = '  -3.44%  '
$ws.Range('E29').Value = '  -1.27%  '
$ws.Range('E30').Value = '  +0.11%  '
$ws.Range('E31').Value = '  -3.27%  '
$ws.Range('E32').Value = '  -6.83%  '
$ws.Range('E33').Value = '  -3.35%  '
$ws.Range('E34').Value = '  +1.76%  '
$ws.Range('E35').Value = '  -15.20%  '
$ws.Range('E36').Value = '  -3.58%  '
$ws.Range('E37').Value = '  -4.03%  '
$ws.Range('E38').Value = '  -2.98%  '
$ws.Range('E39').Value = '  -3.77%  '
$ws.Range('E40').Value = '  -2.47%  '
$ws.Range('E41').Value = '  -4.43%  '
$ws.Range('E42').Value = '  +0.05%  '
$ws.Range('E43').Value = '  -2.68%  '
$ws.Range('E44').Value = '  -3.70%  '
$ws.Range('B45').Value = 'PancakeSwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('E45').Value = '  +0.42%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('E46').Value = '  -5.51%  '
$ws.Range('E47').Value = '  -5.45%  '
$ws.Range('E48').Value = '  -3.36%  '
$ws.Range('E49').Value = '  -4.10%  '
$ws.Range('E50').Value = '  -3.23%  '
$ws.Range('E51').Value = '  -3.58%  '
